$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be auto-parsed as a number
$textCells = @('D5', 'D6', 'D8', 'D9', 'D11', 'D12', 'D14', 'D20', 'D21', 'D22', 'D24', 'D28', 'D30', 'D31', 'D32', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D48', 'D49', 'D50', 'D51')
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range('D2').Value = '66.991.09'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '3.556.71'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '609.07'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = '146.02'
$ws.Range('D7').Value = '3.557.17'
$ws.Range('E7').Value = '  -1.09%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  +5.65%  '
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('D11').Value = '7.88'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = '0.415'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '4.160.08'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('D14').Value = '0.0000198'
$ws.Range('E14').Value = '  -5.38%  '
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('D16').Value = '3.558.10'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('D18').Value = '66.790.41'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  -3.90%  '
$ws.Range('D20').Value = '6.25'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('D21').Value = '14.80'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').Value = '428.21'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('E23').Value = '  -2.71%  '
$ws.Range('D24').Value = '77.81'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').Value = '3.701.33'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D28').Value = '8.12'
$ws.Range('E28').Value = '  -1.79%  '
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('D30').Value = '9.10'
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.158'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('B33').Value = 'RenzoRestakedETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D33').Value = '3.564.48'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').Value = '24.54'
$ws.Range('E34').Value = '  -3.57%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '1.37'
$ws.Range('E36').Value = '  -6.76%  '
$ws.Range('D37').Value = '7.74'
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('D38').Value = '1.66'
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('D39').Value = '177.77'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = '5.33'
$ws.Range('E40').Value = '  -5.50%  '
$ws.Range('D41').Value = '0.0833'
$ws.Range('E41').Value = '  -2.78%  '
$ws.Range('D42').Value = '5.05'
$ws.Range('E42').Value = '  -3.28%  '
$ws.Range('E43').Value = '  -3.50%  '
$ws.Range('D44').Value = '45.68'
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('E45').Value = '  -5.56%  '
$ws.Range('E47').Value = '  -5.42%  '
$ws.Range('D48').Value = '23.73'
$ws.Range('E48').Value = '  -2.57%  '
$ws.Range('D49').Value = '7.17'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').Value = '1.13'
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('D51').Value = '0.927'
$ws.Range('E51').Value = '  -2.64%  '
